# Insert two new data rows at position 320 (pushing the existing rows
# 320-384 down to 322-386), then populate the two new rows with the
# new "Tuna" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 320.
$ws.Rows.Item(320).Insert()
$ws.Rows.Item(320).Insert()

# --- New row 320 ---------------------------------------------------
$ws.Cells.Item(320, 1).Value = 10
$ws.Cells.Item(320, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(320, 3).Value = "La Araucanía"
$ws.Cells.Item(320, 4).Value = 44511
$ws.Cells.Item(320, 5).Value = 9
$ws.Cells.Item(320, 6).Value = 100112027
$ws.Cells.Item(320, 7).Value = "Melón"
$ws.Cells.Item(320, 8).Value = "Tuna"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 40
$ws.Cells.Item(320, 11).Value = 25000
$ws.Cells.Item(320, 12).Value = 25000
$ws.Cells.Item(320, 13).Value = 25000
$ws.Cells.Item(320, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(320, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(320, 16).Value = 1562
$ws.Cells.Item(320, 17).Value = 16
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# --- New row 321 ---------------------------------------------------
$ws.Cells.Item(321, 1).Value = 10
$ws.Cells.Item(321, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(321, 3).Value = "La Araucanía"
$ws.Cells.Item(321, 4).Value = 44511
$ws.Cells.Item(321, 5).Value = 9
$ws.Cells.Item(321, 6).Value = 100112027
$ws.Cells.Item(321, 7).Value = "Melón"
$ws.Cells.Item(321, 8).Value = "Tuna"
$ws.Cells.Item(321, 9).Value = "Segunda"
$ws.Cells.Item(321, 10).Value = 30
$ws.Cells.Item(321, 11).Value = 25000
$ws.Cells.Item(321, 12).Value = 25000
$ws.Cells.Item(321, 13).Value = 25000
$ws.Cells.Item(321, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(321, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(321, 16).Value = 1042
$ws.Cells.Item(321, 17).Value = 24
$ws.Cells.Item(321, 18).Value = "Hortaliza"
